$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($cellRef, $val) {
    $rng = $ws.Range($cellRef)
    $rng.NumberFormat = "@"
    $rng.Value = $val
    $rng.ClearFormats()
}

Set-TextValue "D2" "30.340.64"
Set-TextValue "E2" "  -0.63%  "
Set-TextValue "D3" "1.867.72"
Set-TextValue "E3" "  -0.29%  "
Set-TextValue "E4" "  -0.05%  "
Set-TextValue "D5" "243.66"
Set-TextValue "E5" "  -1.60%  "
Set-TextValue "E6" "  -0.04%  "
Set-TextValue "D7" "0.4703"
Set-TextValue "E7" "  -1.14%  "
Set-TextValue "D8" "0.2871"
Set-TextValue "E8" "  -1.28%  "
Set-TextValue "D9" "0.06446"
Set-TextValue "E9" "  -0.73%  "
Set-TextValue "D10" "22.00"
Set-TextValue "E10" "  +0.78%  "
Set-TextValue "D11" "0.07767"
Set-TextValue "E11" "  +0.28%  "
Set-TextValue "D12" "1.873.83"
Set-TextValue "E12" "  +0.07%  "
Set-TextValue "D13" "95.86"
Set-TextValue "E13" "  +0.01%  "
Set-TextValue "D14" "0.7202"
Set-TextValue "E14" "  -2.16%  "
Set-TextValue "D15" "5.120"
Set-TextValue "E15" "  -1.09%  "
Set-TextValue "D16" "278.17"
Set-TextValue "E16" "  +1.55%  "
Set-TextValue "D17" "30.332.41"
Set-TextValue "E17" "  -0.83%  "
Set-TextValue "D18" "12.95"
Set-TextValue "E18" "  -1.79%  "
Set-TextValue "D19" "0.000007513"
Set-TextValue "E19" "  +0.27%  "
Set-TextValue "D20" "1.000"
Set-TextValue "D21" "2.114.08"
Set-TextValue "E21" "  -0.23%  "
Set-TextValue "E22" "  -0.11%  "
Set-TextValue "D23" "5.212"
Set-TextValue "E23" "  +0.16%  "
Set-TextValue "D24" "6.212"
Set-TextValue "E24" "  +0.82%  "
Set-TextValue "D25" "163.28"
Set-TextValue "E25" "  -0.92%  "
Set-TextValue "D26" "9.026"
Set-TextValue "E26" "  -1.54%  "
Set-TextValue "D27" "18.63"
Set-TextValue "E27" "  -0.65%  "
Set-TextValue "D28" "1.870"
Set-TextValue "E28" "  -1.75%  "
Set-TextValue "E29" "  -1.50%  "
Set-TextValue "D30" "0.09591"
Set-TextValue "D31" "1.466"
Set-TextValue "E31" "  -2.02%  "
Set-TextValue "D32" "4.196"
Set-TextValue "E32" "  -1.32%  "
Set-TextValue "D33" "4.091"
Set-TextValue "E33" "  +0.20%  "
Set-TextValue "D34" "0.04804"
Set-TextValue "E34" "  +0.34%  "
Set-TextValue "D35" "1.116"
Set-TextValue "E35" "  -0.27%  "
Set-TextValue "D36" "0.6883"
Set-TextValue "E37" "  -0.18%  "
Set-TextValue "D38" "0.01872"
Set-TextValue "E38" "  +1.04%  "
Set-TextValue "D39" "2.804"
Set-TextValue "E39" "  +1.83%  "
Set-TextValue "D41" "74.06"
Set-TextValue "E41" "  +0.99%  "
Set-TextValue "D42" "1.934"
Set-TextValue "E42" "  -2.41%  "
Set-TextValue "D43" "0.4208"
Set-TextValue "E43" "  +0.71%  "
Set-TextValue "D44" "0.9992"
Set-TextValue "E44" "  -0.13%  "
Set-TextValue "D45" "0.8280"
Set-TextValue "E45" "  -0.71%  "
Set-TextValue "D46" "100.69"
Set-TextValue "E46" "  -0.84%  "
Set-TextValue "D47" "9.582"
Set-TextValue "E47" "  +1.78%  "
Set-TextValue "D48" "35.22"
Set-TextValue "D49" "6.938"
Set-TextValue "E49" "  -0.20%  "
Set-TextValue "D50" "899.57"
Set-TextValue "E50" "  -1.60%  "
Set-TextValue "D51" "0.05714"
Set-TextValue "E51" "  +0.81%  "
